# Update cryptocurrency Price (D) and Volume(1h) (E) columns
# to match the latest scrape, per GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "29.545.05"
$ws.Range('E2').Value = "  +0.37%  "
$ws.Range('D3').Value = "1.926.98"
$ws.Range('E3').Value = "  +0.96%  "
$ws.Range('E4').Value = "  +0.55%  "
$ws.Range('D5').Value = "'326.57"
$ws.Range('E5').Value = "  +0.44%  "
$ws.Range('E6').Value = "  +0.56%  "
$ws.Range('D7').Value = "'0.4814"
$ws.Range('E7').Value = "  -0.15%  "
$ws.Range('D8').Value = "'0.4049"
$ws.Range('E8').Value = "  -0.25%  "
$ws.Range('D9').Value = "'0.08182"
$ws.Range('E9').Value = "  +0.34%  "
$ws.Range('D10').Value = "'1.008"
$ws.Range('E10').Value = "  -0.32%  "
$ws.Range('D11').Value = "'23.84"
$ws.Range('E11').Value = "  +1.83%  "
$ws.Range('D12').Value = "1.934.42"
$ws.Range('E12').Value = "  +1.01%  "
$ws.Range('D13').Value = "'6.096"
$ws.Range('E13').Value = "  +1.40%  "
$ws.Range('D14').Value = "'7.306"
$ws.Range('E14').Value = "  +2.09%  "
$ws.Range('E15').Value = "  +1.47%  "
$ws.Range('D16').Value = "'0.06896"
$ws.Range('E16').Value = "  +1.55%  "
$ws.Range('D17').Value = "'1.013"
$ws.Range('E17').Value = "  +0.55%  "
$ws.Range('D18').Value = "'0.00001038"
$ws.Range('E18').Value = "  +0.23%  "
$ws.Range('D19').Value = "'17.63"
$ws.Range('E19').Value = "  -0.16%  "
$ws.Range('E20').Value = "  +0.44%  "
$ws.Range('D21').Value = "29.545.41"
$ws.Range('E21').Value = "  +0.30%  "
$ws.Range('D22').Value = "'5.675"
$ws.Range('E22').Value = "  +0.88%  "
$ws.Range('E23').Value = "  +2.58%  "
$ws.Range('D24').Value = "'2.181"
$ws.Range('E24').Value = "  +0.25%  "
$ws.Range('D25').Value = "2.145.61"
$ws.Range('E25').Value = "  -0.12%  "
$ws.Range('D26').Value = "'156.15"
$ws.Range('E26').Value = "  -0.10%  "
$ws.Range('D27').Value = "'6.364"
$ws.Range('E27').Value = "  -1.94%  "
$ws.Range('D28').Value = "'20.02"
$ws.Range('E28').Value = "  -0.12%  "
$ws.Range('D29').Value = "'2.090"
$ws.Range('E29').Value = "  -0.91%  "
$ws.Range('D30').Value = "'120.55"
$ws.Range('E30').Value = "  +0.14%  "
$ws.Range('D31').Value = "'1.012"
$ws.Range('E31').Value = "  -1.50%  "
$ws.Range('D32').Value = "'0.09577"
$ws.Range('D33').Value = "'5.592"
$ws.Range('E33').Value = "  +1.49%  "
$ws.Range('D34').Value = "'3.562"
$ws.Range('E34').Value = "  -0.03%  "
$ws.Range('D35').Value = "'1.383"
$ws.Range('E35').Value = "  -0.42%  "
$ws.Range('D36').Value = "'0.06337"
$ws.Range('E36').Value = "  +3.85%  "
$ws.Range('D37').Value = "'0.02281"
$ws.Range('E37').Value = "  +0.56%  "
$ws.Range('D38').Value = "'1.192"
$ws.Range('E38').Value = "  +1.42%  "
$ws.Range('D39').Value = "'0.5942"
$ws.Range('E39').Value = "  -0.05%  "
$ws.Range('D40').Value = "'10.71"
$ws.Range('E40').Value = "  -0.88%  "
$ws.Range('D41').Value = "'1.011"
$ws.Range('E41').Value = "  +0.55%  "
$ws.Range('D42').Value = "'7.866"
$ws.Range('E42').Value = "  -1.31%  "
$ws.Range('D43').Value = "'0.1843"
$ws.Range('E43').Value = "  -0.59%  "
$ws.Range('D44').Value = "'2.460"
$ws.Range('E44').Value = "  +3.53%  "
$ws.Range('D45').Value = "'1.246"
$ws.Range('E45').Value = "  -2.25%  "
$ws.Range('D46').Value = "'12.37"
$ws.Range('E46').Value = "  -1.52%  "
$ws.Range('D47').Value = "'0.07487"
$ws.Range('E47').Value = "  -1.55%  "
$ws.Range('D48').Value = "'0.5543"
$ws.Range('E48').Value = "  -0.31%  "
$ws.Range('D49').Value = "'1.970"
$ws.Range('E49').Value = "  +1.25%  "
$ws.Range('D50').Value = "'117.67"
$ws.Range('E50').Value = "  +1.01%  "
$ws.Range('D51').Value = "'2.434"
$ws.Range('E51').Value = "  +1.24%  "
